# Apply the "added release table and some updated data" commit to the workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# 1) Small data corrections on existing Sheet1 rows
# ---------------------------------------------------------------------------
$ws1.Range("K17").Value = 3490
$ws1.Range("K22").Value = 3490

$ws1.Range("I45").Value = 815
$ws1.Range("K45").Value = 3307

$ws1.Range("I50").Value = 815
$ws1.Range("K50").Value = 3307

$ws1.Range("I87").Value = 801
$ws1.Range("K87").Value = 3302

$ws1.Range("I92").Value = 801
$ws1.Range("K92").Value = 3302

$ws1.Range("I124").Value = 768
$ws1.Range("K124").Value = 3115

$ws1.Range("I125").Value = 379
$ws1.Range("J125").Value = 34
$ws1.Range("K125").Value = 1596

$ws1.Range("H126").Value = 70931

$ws1.Range("H127").Value = 84693

$ws1.Range("H129").Value = 546427
$ws1.Range("I129").Value = 1147
$ws1.Range("J129").Value = 188
$ws1.Range("K129").Value = 3859

# ---------------------------------------------------------------------------
# 2) Append two new weeks (20 and 21) of survey data on Sheet1, rows 130-143
# ---------------------------------------------------------------------------
$week20 = @(
    @(20, "Suisun Bay",          6,  24, 0, 0, 0, 115859, $null, $null, $null),
    @(20, "Suisun Marsh",        6,  22, 1, 0, 0, 102403, 817,   82,    3336),
    @(20, "Lower Sacramento",    6,  24, 0, 0, 0, 94714,  $null, $null, $null),
    @(20, "Cache Slough LI",     6,  23, 0, 0, 0, 84701,  $null, $null, $null),
    @(20, "Sac DW Ship Channel", 5,  20, 0, 0, 0, 72976,  $null, $null, $null),
    @(20, "Lower San Joaquin",   6,  24, 0, 0, 0, 95140,  $null, $null, $null),
    @(20, "All Strata",          35, 137, 1, 0, 0, 565794, 817,   82,    3336)
)

$week21 = @(
    @(21, "Suisun Bay",          6,  24, 0, 0, 0, 111321, $null, $null, $null),
    @(21, "Suisun Marsh",        6,  22, 0, 0, 0, 115261, $null, $null, $null),
    @(21, "Lower Sacramento",    5,  19, 0, 0, 0, 77035,  $null, $null, $null),
    @(21, "Cache Slough LI",     6,  22, 0, 0, 0, 71683,  $null, $null, $null),
    @(21, "Sac DW Ship Channel", 6,  24, 0, 0, 0, 92880,  $null, $null, $null),
    @(21, "Lower San Joaquin",   6,  24, 0, 0, 0, 94594,  $null, $null, $null),
    @(21, "All Strata",          35, 135, 0, 0, 0, 562775, $null, $null, $null)
)

$newWeeks = @($week20, $week21)

$r = 130
foreach ($week in $newWeeks) {
    foreach ($row in $week) {
        $ws1.Cells.Item($r, 1).Value = $row[0]
        $ws1.Cells.Item($r, 2).Value = $row[1]
        $ws1.Cells.Item($r, 3).Value = $row[2]
        $ws1.Cells.Item($r, 4).Value = $row[3]
        $ws1.Cells.Item($r, 5).Value = $row[4]
        $ws1.Cells.Item($r, 6).Value = $row[5]
        $ws1.Cells.Item($r, 7).Value = $row[6]

        $hCell = $ws1.Cells.Item($r, 8)
        $hCell.Value = $row[7]
        $hCell.NumberFormat = "#,##0"

        if ($null -eq $row[8]) {
            $ws1.Cells.Item($r, 9).Value = "0*"
            $ws1.Cells.Item($r, 10).Value = "NA"
            $ws1.Cells.Item($r, 11).Value = "NA"
        } else {
            $ws1.Cells.Item($r, 9).Value = $row[8]
            $ws1.Cells.Item($r, 10).Value = $row[9]
            $kCell = $ws1.Cells.Item($r, 11)
            $kCell.Value = $row[10]
            $kCell.NumberFormat = "#,##0"
        }

        $r = $r + 1
    }
}

# Row 129's "All Strata" abundance index (I129) is >= 1000, so it keeps the
# thousands-separator number format just like K129/H129.
$ws1.Range("I129").NumberFormat = "#,##0"

# ---------------------------------------------------------------------------
# 3) Append the matching release-date rows on Sheet2 (rows 21-22)
# ---------------------------------------------------------------------------
$ws2.Cells.Item(21, 1).Value = 20
$ws2.Cells.Item(21, 2).Value = "October 13–17, 2025"

$ws2.Cells.Item(22, 1).Value = 21
$ws2.Cells.Item(22, 2).Value = "October 20–24, 2025"

# ---------------------------------------------------------------------------
# 4) Update view/selection state to match the saved workbook
# ---------------------------------------------------------------------------
[void]$ws1.Range("A1:K143").Select()
[void]$ws2.Range("I9").Select()
[void]$ws2.Activate()
